$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 15 height change (128 -> 112) ---
$ws.Rows.Item(15).RowHeight = 112

# --- G58: remove bottom border (visual tweak that accompanies the new row 59 data) ---
$ws.Range("G58").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom = 9, xlLineStyleNone = -4142

# --- Row 59: fill in the new literature entry (Heffel/Finnigan follow-up: Cleave & Rescue) ---
$ws.Range("A59").Value2 = 57
$ws.Range("B59").Value2 = "Peer reviewed"
$ws.Range("C59").Value2 = "PNAS"

# D59 needs to be the literal text "0.5" (not the number 0.5), matching the rest of column D.
$ws.Range("D7").Copy()
$ws.Range("D59").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("E59").Value2 = 1
$ws.Range("F59").Value2 = 2020

# H59: DOI text + live hyperlink, restyled to match the other DOI cells (style of H56).
$ws.Range("H59").Value2 = "https://doi.org/10.1073/pnas.1921698117"
$ws.Hyperlinks.Add($ws.Range("H59"), "https://doi.org/10.1073/pnas.1921698117") | Out-Null
$ws.Range("H56").Copy()
$ws.Range("H59").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I59").Value2 = 10

# J59 picks up the "left/center, no-wrap" style already used on J4.
$ws.Range("J59").Value2 = "Fruit fly"
$ws.Range("J4").Copy()
$ws.Range("J59").PasteSpecial(-4122)

# K59 picks up the wrap style already used on K3.
$ws.Range("K59").Value2 = "Cleave and Rescue, CRIPR"
$ws.Range("K3").Copy()
$ws.Range("K59").PasteSpecial(-4122)

# L59 picks up the style already used on L4.
$ws.Range("L59").Value2 = "Experiment"
$ws.Range("L4").Copy()
$ws.Range("L59").PasteSpecial(-4122)

$ws.Range("M59").Value2 = "NA"
$ws.Range("N59").Value2 = "NA"
$ws.Range("O59").Value2 = "NA"
$ws.Range("P59").Value2 = "Replacement"
$ws.Range("Q59").Value2 = 1

$excel.CutCopyMode = 0

# --- Sheet view: scroll down to the newly-edited row and select R59 ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 57 } catch {}
$ws.Range("R59").Select()
